$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 95.63567833333333
$ws.Range("H2").Value = 286.907035
$ws.Range("I2").Value = 0.2808828217467972
$ws.Range("J2").Value = 0.2808828217467972
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.372179
$ws.Range("N2").Value = 4.116537
$ws.Range("O2").Value = 0.2533567233062949
$ws.Range("P2").Value = 0.2533567233062949
$ws.Range("Q2").Value = 131.229269459755
$ws.Range("R2").Value = 1181.063425137795
$ws.Range("S2").Value = 0.07116355135079466
$ws.Range("T2").Value = 0.07116355135079468

$ws.Range("G3").Value = 95.63567833333333
$ws.Range("H3").Value = 286.907035
$ws.Range("I3").Value = 0.2808828217467972
$ws.Range("J3").Value = 0.2808828217467972
$ws.Range("O3").Value = 0.01171683533985869
$ws.Range("P3").Value = 0.0117168353398587
$ws.Range("Q3").Value = 6.06888075423611
$ws.Range("R3").Value = 54.619926788125
$ws.Range("S3").Value = 0.003291057772202103
$ws.Range("T3").Value = 0.003291057772202104

$ws.Range("G4").Value = 95.63567833333333
$ws.Range("H4").Value = 286.907035
$ws.Range("I4").Value = 0.2808828217467972
$ws.Range("J4").Value = 0.2808828217467972
$ws.Range("M4").Value = 3.980358666666666
$ws.Range("O4").Value = 0.7349264413538463
$ws.Range("P4").Value = 0.7349264413538463
$ws.Range("Q4").Value = 380.6643010966288
$ws.Range("R4").Value = 3425.97870986966
$ws.Range("S4").Value = 0.2064282126238004
$ws.Range("T4").Value = 0.2064282126238005

$ws.Range("I5").Value = 0.392628215788982
$ws.Range("J5").Value = 0.392628215788982
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.372179
$ws.Range("N5").Value = 4.116537
$ws.Range("O5").Value = 0.2533567233062949
$ws.Range("P5").Value = 0.2533567233062949
$ws.Range("Q5").Value = 183.43704184844
$ws.Range("R5").Value = 1650.93337663596
$ws.Range("S5").Value = 0.09947499822989338
$ws.Range("T5").Value = 0.09947499822989338

$ws.Range("I6").Value = 0.392628215788982
$ws.Range("J6").Value = 0.392628215788982
$ws.Range("O6").Value = 0.01171683533985869
$ws.Range("P6").Value = 0.0117168353398587
$ws.Range("S6").Value = 0.00460036015418201
$ws.Range("T6").Value = 0.00460036015418201

$ws.Range("I7").Value = 0.392628215788982
$ws.Range("J7").Value = 0.392628215788982
$ws.Range("M7").Value = 3.980358666666666
$ws.Range("O7").Value = 0.7349264413538463
$ws.Range("P7").Value = 0.7349264413538463
$ws.Range("Q7").Value = 532.1063937788978
$ws.Range("R7").Value = 4788.95754401008
$ws.Range("S7").Value = 0.2885528574049066
$ws.Range("T7").Value = 0.2885528574049066

$ws.Range("I8").Value = 0.3264889624642208
$ws.Range("J8").Value = 0.3264889624642208
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.372179
$ws.Range("N8").Value = 4.116537
$ws.Range("O8").Value = 0.2533567233062949
$ws.Range("P8").Value = 0.2533567233062949
$ws.Range("Q8").Value = 152.53659075483
$ws.Range("R8").Value = 1372.82931679347
$ws.Range("S8").Value = 0.0827181737256069
$ws.Range("T8").Value = 0.08271817372560691

$ws.Range("I9").Value = 0.3264889624642208
$ws.Range("J9").Value = 0.3264889624642208
$ws.Range("O9").Value = 0.01171683533985869
$ws.Range("P9").Value = 0.0117168353398587
$ws.Range("S9").Value = 0.00382541741347458
$ws.Range("T9").Value = 0.003825417413474582

$ws.Range("I10").Value = 0.3264889624642208
$ws.Range("J10").Value = 0.3264889624642208
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.980358666666666
$ws.Range("O10").Value = 0.7349264413538463
$ws.Range("P10").Value = 0.7349264413538463
$ws.Range("Q10").Value = 442.4716753388399
$ws.Range("S10").Value = 0.2399453713251393
$ws.Range("T10").Value = 0.2399453713251393

